$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G, shifting "Created At"/"Updated At" right.
$ws.Columns("G:G").Insert()

# New header for the inserted column.
$ws.Range("G1").Value = "Merge Requests"

# New values for the inserted column's data rows.
$ws.Range("G2").Value = 'Resolve "Can we post the log on a port through a socket ?"'
$ws.Range("G3").Value = 'Draft: Resolve "Problem after 24H"; Resolve "Problem after 24H"'
